$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) dates in rows 2-6 from 45183 to 45184
$ws.Range("C2:C6").Value = 45184
